# Auto-generated Excel COM-interop script
# Applies updated crypto price/volume(1h) values to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.983.75'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.459.86'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  -0.13%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '520.14'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -2.45%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '133.00'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.32%  '
$ws.Range('E7').Value = '  +0.02%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.556'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').Value = '2.467.90'
$ws.Range('E9').Value = '  -1.66%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0976'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -3.38%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.156'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('E12').Value = '  -2.33%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.338'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -2.94%  '
$ws.Range('D14').Value = '2.901.52'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').Value = '57.911.97'
$ws.Range('E15').Value = '  -1.39%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '22.25'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -2.84%  '
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').Value = '2.474.49'
$ws.Range('E18').Value = '  -1.64%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.59'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -4.15%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '319.57'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.90%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '4.14'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -2.40%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  -4.53%  '
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('E25').Value = '  -2.88%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('E27').Value = '  -3.16%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.31'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('D29').Value = '0.0₃0744'
$ws.Range('E29').Value = '  -2.60%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '167.37'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('E31').Value = '  -3.17%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '6.20'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -5.18%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.16'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  -0.07%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +0.05%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '1.35'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -0.75%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '17.99'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -1.83%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '3.95'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -2.06%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '36.29'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('E40').Value = '  -4.48%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.792'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -1.79%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '272.00'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -3.96%  '
$ws.Range('E44').Value = '  -3.38%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.587'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -3.19%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '123.65'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -5.05%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.0906'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -1.67%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.0486'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -3.07%  '
$ws.Range('E49').Value = '  -2.94%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '16.80'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('D51').Value = '1.720.98'
$ws.Range('E51').Value = '  -1.89%  '
